# Auto-update price data: insert a new "today" row at the top of the
# price table (row 2), pushing all existing date rows down by one.
# The new row carries the same commodity values as the rest of the
# series (unchanged market prices), dated one day after the previous
# top row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record just below the header row.
$ws.Rows.Item(2).Insert()

# Column A holds the date as plain text (e.g. "2026-01-23"), not a
# real Excel date. Force the new cell to be treated as text too, so
# Excel doesn't silently convert the "2026-01-24" string into a date
# serial number; then drop back to the workbook's default style so no
# stray per-cell formatting is introduced.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value2 = "2026-01-24"
$ws.Range("A2").Style = "Normal"

# Same commodity prices as every other row in this (currently flat)
# series.
$ws.Range("B2").Value2 = 783.5
$ws.Range("C2").Value2 = 1112
$ws.Range("D2").Value2 = 3610
